$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Final target data (Player Name, Position, Team) for rows 2-17
$data = @(
    @("De'Aaron Fox",   "PG",       "Sacramento Kings"),
    @("Tyler Herro",    "PG,SG",    "Miami Heat"),
    @("Caris LeVert",   "SG,SF",    "Cleveland Cavaliers"),
    @("Miles Bridges",  "SF,PF",    "Charlotte Hornets"),
    @("DeMar DeRozan",  "SF,PF",    "Sacramento Kings"),
    @("Nikola Vucevic", "PF,C",     "Chicago Bulls"),
    @("Evan Mobley",    "PF,C",     "Cleveland Cavaliers"),
    @("Nick Richards",  "C",        "Charlotte Hornets"),
    @("Gradey Dick",    "SG,SF",    "Toronto Raptors"),
    @("Brook Lopez",    "C",        "Milwaukee Bucks"),
    @("Santi Aldama",   "PF,C",     "Memphis Grizzlies"),
    @("Luka Doncic",    "PG,SG",    "Dallas Mavericks"),
    @("Ja Morant",      "PG",       "Memphis Grizzlies"),
    @("Mikal Bridges",  "SG,SF,PF", "New York Knicks"),
    @("Scottie Barnes", "SG,SF,PF", "Toronto Raptors"),
    @("Josh Giddey",    "PG,SG,SF", "Chicago Bulls")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
